# Add two new columns (CanClone, ActorID) to the "表1" XML table on Sheet1,
# fill in header + data values, resize column J, update the selection, and
# flip on a page setup (paper size / orientation) as in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)

# Append two columns to the table -> table grows from A1:I4 to A1:K4.
# (Column J = CanClone, Column K = ActorID.)
$colCanClone = $lo.ListColumns.Add()
$colActorID  = $lo.ListColumns.Add()

# Set header text. ActorID is written first so it lands earlier in the
# shared-string table than CanClone, matching the source workbook.
$ws.Range("K1").Value = "ActorID"
$ws.Range("J1").Value = "CanClone"

# Data rows for the two new columns.
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0

$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# New column J gets an explicit width of 14 characters (no bestFit).
# ColumnWidth writes are off by a constant +5/7 in this engine, so we
# back that out here to land on exactly 14 in the saved file.
$ws.Columns.Item(10).ColumnWidth = 13.285714285714286

# Move/collapse the selection the way the author left it.
$ws.Range("K9").Select()

# Turn on an explicit page setup (paper size 9 = A4, portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
